$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2014.6666
$ws.Range("I12").Value = 2014.6666
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 2014.6666
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -1844.6666
$ws.Range("N12").Value = $null
$ws.Range("H58").Value = 3922.6667
$ws.Range("I58").Value = 663
$ws.Range("K58").Value = 1989
$ws.Range("M58").Value = -1839
$ws.Range("H62").Value = 1496.1818
$ws.Range("I62").Value = 1555.8
$ws.Range("J62").Value = 900
$ws.Range("K62").Value = 1555.8
$ws.Range("L62").Value = 900
$ws.Range("M62").Value = -931.8
$ws.Range("N62").Value = -2148
$ws.Range("H65").Value = 1496.1818
$ws.Range("I65").Value = 1555.8
$ws.Range("J65").Value = 900
$ws.Range("K65").Value = 7779
$ws.Range("L65").Value = 4500
$ws.Range("M65").Value = -4659
$ws.Range("N65").Value = -10740
$ws.Range("H70").Value = 2501
$ws.Range("I70").Value = 1583.3334
$ws.Range("J70").Value = 3877.5
$ws.Range("K70").Value = 4750.0002
$ws.Range("L70").Value = 11632.5
$ws.Range("M70").Value = -4480.0002
$ws.Range("N70").Value = -12172.5
$ws.Range("H73").Value = 2501
$ws.Range("I73").Value = 1583.3334
$ws.Range("J73").Value = 3877.5
$ws.Range("K73").Value = 4750.0002
$ws.Range("L73").Value = 11632.5
$ws.Range("M73").Value = -3814.0002
$ws.Range("N73").Value = -13504.5
$ws.Range("H129").Value = 862.24
$ws.Range("I129").Value = 500
$ws.Range("J129").Value = 865.899
$ws.Range("K129").Value = 1500
$ws.Range("L129").Value = 2597.697
$ws.Range("M129").Value = 3500
$ws.Range("N129").Value = -12597.697
$ws.Range("H132").Value = 32586034
$ws.Range("I132").Value = 37042450
$ws.Range("J132").Value = 2505227.8
$ws.Range("K132").Value = 111127350
$ws.Range("L132").Value = 7515683.399999999
$ws.Range("M132").Value = -111124820
$ws.Range("N132").Value = -7520743.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2104.5454
$ws.Range("I45").Value = 2308.875
$ws.Range("J45").Value = 1559.6666
$ws.Range("K45").Value = 2308.875
$ws.Range("L45").Value = 1559.6666
$ws.Range("M45").Value = -1931.875
$ws.Range("N45").Value = -2313.6666
$ws.Range("H74").Value = 4114.6
$ws.Range("I74").Value = 3656.2727
$ws.Range("K74").Value = 3656.2727
$ws.Range("M74").Value = -2782.2727
$ws.Range("H77").Value = 4114.6
$ws.Range("I77").Value = 3656.2727
$ws.Range("K77").Value = 18281.3635
$ws.Range("M77").Value = -13913.3635
$ws.Range("H122").Value = 4461.079
$ws.Range("I122").Value = 3848.9312
$ws.Range("J122").Value = 6433.5557
$ws.Range("K122").Value = 11546.7936
$ws.Range("L122").Value = 19300.6671
$ws.Range("M122").Value = -9096.793600000001
$ws.Range("N122").Value = -24200.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 954.73334
$ws.Range("I64").Value = 819.125
$ws.Range("J64").Value = 1109.7142
$ws.Range("K64").Value = 819.125
$ws.Range("L64").Value = 1109.7142
$ws.Range("M64").Value = -594.125
$ws.Range("N64").Value = -1559.7142
$ws.Range("H67").Value = 954.73334
$ws.Range("I67").Value = 819.125
$ws.Range("J67").Value = 1109.7142
$ws.Range("K67").Value = 819.125
$ws.Range("L67").Value = 1109.7142
$ws.Range("M67").Value = -39.125
$ws.Range("N67").Value = -2669.7142
$ws.Range("H94").Value = 1178.5714
$ws.Range("I94").Value = 1380
$ws.Range("J94").Value = 675
$ws.Range("K94").Value = 1380
$ws.Range("L94").Value = 675
$ws.Range("M94").Value = -929
$ws.Range("N94").Value = -1577

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4467.185
$ws.Range("I31").Value = 1792.4878
$ws.Range("J31").Value = 12902.77
$ws.Range("K31").Value = 1792.4878
$ws.Range("L31").Value = 12902.77
$ws.Range("M31").Value = -1497.4878
$ws.Range("N31").Value = -13492.77
$ws.Range("H34").Value = 4467.185
$ws.Range("I34").Value = 1792.4878
$ws.Range("J34").Value = 12902.77
$ws.Range("K34").Value = 1792.4878
$ws.Range("L34").Value = 12902.77
$ws.Range("M34").Value = -1590.4878
$ws.Range("N34").Value = -13306.77
$ws.Range("H58").Value = 1765.5405
$ws.Range("I58").Value = 1349.1936
$ws.Range("J58").Value = 3916.6667
$ws.Range("K58").Value = 1349.1936
$ws.Range("L58").Value = 3916.6667
$ws.Range("M58").Value = -1146.1936
$ws.Range("N58").Value = -4322.6667
$ws.Range("H136").Value = 1765.5405
$ws.Range("I136").Value = 1349.1936
$ws.Range("J136").Value = 3916.6667
$ws.Range("K136").Value = 4047.5808
$ws.Range("L136").Value = 11750.0001
$ws.Range("M136").Value = -1497.5808
$ws.Range("N136").Value = -16850.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 557503.4399999999
$ws.Range("J5").Value = 1028855.94
$ws.Range("L5").Value = 3086567.82
$ws.Range("N5").Value = -3086791.82
$ws.Range("H93").Value = 9727.272000000001
$ws.Range("J93").Value = 9727.272000000001
$ws.Range("L93").Value = 29181.816
$ws.Range("N93").Value = -32925.81600000001
$ws.Range("H113").Value = 6250598.5
$ws.Range("I113").Value = 624.6667
$ws.Range("J113").Value = 15625560
$ws.Range("K113").Value = 1874.0001
$ws.Range("L113").Value = 46876680
$ws.Range("M113").Value = 295.9999
$ws.Range("N113").Value = -46881020
$ws.Range("H122").Value = 3313.5757
$ws.Range("I122").Value = 1275
$ws.Range("J122").Value = 3766.5925
$ws.Range("K122").Value = 11475
$ws.Range("L122").Value = 33899.3325
$ws.Range("M122").Value = -9025
$ws.Range("N122").Value = -38799.3325
$ws.Range("H131").Value = 780.79
$ws.Range("I131").Value = 315
$ws.Range("J131").Value = 810.5213
$ws.Range("K131").Value = 945
$ws.Range("L131").Value = 2431.5639
$ws.Range("M131").Value = 4095
$ws.Range("N131").Value = -12511.5639
$ws.Range("H132").Value = 2047.3334
$ws.Range("J132").Value = 2789.318
$ws.Range("L132").Value = 25103.862
$ws.Range("N132").Value = -30163.862
$ws.Range("H135").Value = 557503.4399999999
$ws.Range("J135").Value = 1028855.94
$ws.Range("L135").Value = 9259703.459999999
$ws.Range("N135").Value = -9264773.459999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2559.7354
$ws.Range("I132").Value = 1865.6786
$ws.Range("J132").Value = 5798.6665
$ws.Range("K132").Value = 5597.0358
$ws.Range("L132").Value = 17395.9995
$ws.Range("M132").Value = -3067.0358
$ws.Range("N132").Value = -22455.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 4418
$ws.Range("I17").Value = 1900
$ws.Range("J17").Value = 4921.6
$ws.Range("K17").Value = 1900
$ws.Range("L17").Value = 4921.6
$ws.Range("M17").Value = -1730
$ws.Range("N17").Value = -5261.6
$ws.Range("H132").Value = 7339.25
$ws.Range("I132").Value = 3326.2856
$ws.Range("K132").Value = 9978.856800000001
$ws.Range("M132").Value = -7448.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5387.1562
$ws.Range("I122").Value = 3871
$ws.Range("J122").Value = 7603.077
$ws.Range("K122").Value = 11613
$ws.Range("L122").Value = 22809.231
$ws.Range("M122").Value = -9163
$ws.Range("N122").Value = -27709.231
$ws.Range("H126").Value = 2120.4138
$ws.Range("I126").Value = 1173.6
$ws.Range("J126").Value = 4224.4443
$ws.Range("K126").Value = 3520.8
$ws.Range("L126").Value = 12673.3329
$ws.Range("M126").Value = -1050.8
$ws.Range("N126").Value = -17613.3329
$ws.Range("H136").Value = 7560.3213
$ws.Range("I136").Value = 7472
$ws.Range("K136").Value = 22416
$ws.Range("M136").Value = -19866
$ws.Range("H141").Value = 39715
$ws.Range("J141").Value = 39715
$ws.Range("L141").Value = 39715
$ws.Range("N141").Value = -50075
